$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.788933757379432
$ws.Range("D2").Value = 8.631882375455785
$ws.Range("E2").Value = 13.17182002723238
$ws.Range("F2").Value = 28.11950727904313
$ws.Range("G2").Value = 3.627735677940589
$ws.Range("J2").Value = 9.604334820673522
$ws.Range("M2").Value = 43.45304182274057
$ws.Range("O2").Value = 21.85006251816048
$ws.Range("B3").Value = 7.716409253199856
$ws.Range("D3").Value = 8.685057172837677
$ws.Range("E3").Value = 13.15723053124531
$ws.Range("F3").Value = 28.46207500404317
$ws.Range("G3").Value = 3.631117241975744
$ws.Range("J3").Value = 9.618446930900204
$ws.Range("M3").Value = 41.05712227052909
$ws.Range("O3").Value = 22.03613176624007
$ws.Range("B4").Value = 7.673359404913072
$ws.Range("D4").Value = 8.719861654947088
$ws.Range("E4").Value = 13.15381896006366
$ws.Range("F4").Value = 28.6873875739586
$ws.Range("G4").Value = 3.633290978846356
$ws.Range("J4").Value = 9.630516057771258
$ws.Range("M4").Value = 39.50539998493304
$ws.Range("O4").Value = 22.16101529481026
$ws.Range("B5").Value = 7.656207562642632
$ws.Range("D5").Value = 8.734585950316037
$ws.Range("E5").Value = 13.1538175995189
$ws.Range("F5").Value = 28.78292263680365
$ws.Range("G5").Value = 3.634201404761392
$ws.Range("J5").Value = 9.636284579453356
$ws.Range("M5").Value = 38.85306864490125
$ws.Range("O5").Value = 22.21454518959347
$ws.Range("B6").Value = 7.653383689711964
$ws.Range("D6").Value = 8.737063576120672
$ws.Range("E6").Value = 13.15390106482933
$ws.Range("F6").Value = 28.79900939979943
$ws.Range("G6").Value = 3.634354069838296
$ws.Range("J6").Value = 9.637293617810165
$ws.Range("M6").Value = 38.74355162292209
$ws.Range("O6").Value = 22.22359213398844
$ws.Range("B7").Value = 7.673126480080233
$ws.Range("D7").Value = 8.720058041524
$ws.Range("E7").Value = 13.15381332671391
$ws.Range("F7").Value = 28.68866099845459
$ws.Range("G7").Value = 3.633303157378692
$ws.Range("J7").Value = 9.630590419196634
$ws.Range("M7").Value = 39.49668294228341
$ws.Range("O7").Value = 22.16172657980866
$ws.Range("B8").Value = 7.763631660396017
$ws.Range("D8").Value = 8.649769346509752
$ws.Range("E8").Value = 13.16563526747963
$ws.Range("F8").Value = 28.23448772853322
$ws.Range("G8").Value = 3.628881482353275
$ws.Range("J8").Value = 9.608490727341188
$ws.Range("M8").Value = 42.64380855154906
$ws.Range("O8").Value = 21.91198997598805
$ws.Range("B9").Value = 7.951992115157704
$ws.Range("D9").Value = 8.529066773574199
$ws.Range("E9").Value = 13.23303221020233
$ws.Range("F9").Value = 27.46487557626484
$ws.Range("G9").Value = 3.620978753575079
$ws.Range("J9").Value = 9.592413075554587
$ws.Range("M9").Value = 48.16651018787054
$ws.Range("O9").Value = 21.50827615869622
$ws.Range("B10").Value = 8.095818352167354
$ws.Range("D10").Value = 8.450885549833183
$ws.Range("E10").Value = 13.30970471394637
$ws.Range("F10").Value = 26.97639213712159
$ws.Range("G10").Value = 3.615633777035453
$ws.Range("J10").Value = 9.597544737715955
$ws.Range("M10").Value = 51.81987676883953
$ws.Range("O10").Value = 21.26637137744965
$ws.Range("B11").Value = 8.162169137691409
$ws.Range("D11").Value = 8.417612714937261
$ws.Range("E11").Value = 13.35049797304451
$ws.Range("F11").Value = 26.77161222624883
$ws.Range("G11").Value = 3.613300786764014
$ws.Range("J11").Value = 9.603622036444179
$ws.Range("M11").Value = 53.39310810750335
$ws.Range("O11").Value = 21.16869925211117
$ws.Range("B12").Value = 8.187405294361403
$ws.Range("D12").Value = 8.405344197307068
$ws.Range("E12").Value = 13.36679491359078
$ws.Range("F12").Value = 26.69663933244009
$ws.Range("G12").Value = 3.612431381346222
$ws.Range("J12").Value = 9.606466174103028
$ws.Range("M12").Value = 53.97604554318882
$ws.Range("O12").Value = 21.1335343295642
$ws.Range("B13").Value = 8.18196567269497
$ws.Range("D13").Value = 8.407971686151678
$ws.Range("E13").Value = 13.3632473286908
$ws.Range("F13").Value = 26.71267057887746
$ws.Range("G13").Value = 3.612618000392977
$ws.Range("J13").Value = 9.605829428075809
$ws.Range("M13").Value = 53.8510699507992
$ws.Range("O13").Value = 21.14102603785803
$ws.Range("B14").Value = 8.164243242893802
$ws.Range("D14").Value = 8.416596729753802
$ws.Range("E14").Value = 13.35182171939949
$ws.Range("F14").Value = 26.76539221018722
$ws.Range("G14").Value = 3.613228979372992
$ws.Range("J14").Value = 9.603845122011171
$ws.Range("M14").Value = 53.44132370991339
$ws.Range("O14").Value = 21.16576946492472
$ws.Range("B15").Value = 8.153401481676688
$ws.Range("D15").Value = 8.421922996629117
$ws.Range("E15").Value = 13.34493376869582
$ws.Range("F15").Value = 26.79802279407501
$ws.Range("G15").Value = 3.613605047676796
$ws.Range("J15").Value = 9.602700494940827
$ws.Range("M15").Value = 53.18867267416623
$ws.Range("O15").Value = 21.18116399343237
$ws.Range("B16").Value = 8.091498852401187
$ws.Range("D16").Value = 8.453106263577052
$ws.Range("E16").Value = 13.30715773126352
$ws.Range("F16").Value = 26.99013135882948
$ws.Range("G16").Value = 3.615788215286768
$ws.Range("J16").Value = 9.597223287536712
$ws.Range("M16").Value = 51.71527214326752
$ws.Range("O16").Value = 21.27300703889394
$ws.Range("B17").Value = 8.053744281240407
$ws.Range("D17").Value = 8.472824346425737
$ws.Range("E17").Value = 13.28549785573623
$ws.Range("F17").Value = 27.11249484797607
$ws.Range("G17").Value = 3.617152658196565
$ws.Range("J17").Value = 9.594825208039646
$ws.Range("M17").Value = 50.78862203525682
$ws.Range("O17").Value = 21.33254841480112
$ws.Range("B18").Value = 8.032117106385334
$ws.Range("D18").Value = 8.484381316444789
$ws.Range("E18").Value = 13.27359652317716
$ws.Range("F18").Value = 27.18451190603431
$ws.Range("G18").Value = 3.617946724561886
$ws.Range("J18").Value = 9.593797942624875
$ws.Range("M18").Value = 50.24728981490825
$ws.Range("O18").Value = 21.36795694222252
$ws.Range("B19").Value = 8.024810353483465
$ws.Range("D19").Value = 8.488331300420963
$ws.Range("E19").Value = 13.26966258835931
$ws.Range("F19").Value = 27.20917481012355
$ws.Range("G19").Value = 3.618217178150665
$ws.Range("J19").Value = 9.593510447112807
$ws.Range("M19").Value = 50.06257242142582
$ws.Range("O19").Value = 21.38014393612159
$ws.Range("B20").Value = 8.057754347033013
$ws.Range("D20").Value = 8.47070299295161
$ws.Range("E20").Value = 13.28774595785461
$ws.Range("F20").Value = 27.09929917384785
$ws.Range("G20").Value = 3.617006451811716
$ws.Range("J20").Value = 9.595044017956994
$ws.Range("M20").Value = 50.88812994461984
$ws.Range("O20").Value = 21.32608957586957
$ws.Range("B21").Value = 8.169445926521616
$ws.Range("D21").Value = 8.41405434397624
$ws.Range("E21").Value = 13.35515465989478
$ws.Range("F21").Value = 26.74983622374772
$ws.Range("G21").Value = 3.61304913969259
$ws.Range("J21").Value = 9.604413194778486
$ws.Range("M21").Value = 53.56202412981312
$ws.Range("O21").Value = 21.15845194838393
$ws.Range("B22").Value = 8.243077000086291
$ws.Range("D22").Value = 8.37896239026461
$ws.Range("E22").Value = 13.40415979325589
$ws.Range("F22").Value = 26.53647101902636
$ws.Range("G22").Value = 3.610544633427935
$ws.Range("J22").Value = 9.613702190892797
$ws.Range("M22").Value = 55.23488554406673
$ws.Range("O22").Value = 21.05952752466053
$ws.Range("B23").Value = 8.203728002238494
$ws.Range("D23").Value = 8.397514366170283
$ws.Range("E23").Value = 13.37755263012818
$ws.Range("F23").Value = 26.64895034163516
$ws.Range("G23").Value = 3.611873885782747
$ws.Range("J23").Value = 9.608453399844638
$ws.Range("M23").Value = 54.34889387576823
$ws.Range("O23").Value = 21.11133799058097
$ws.Range("B24").Value = 8.055941150305413
$ws.Range("D24").Value = 8.471661369289251
$ws.Range("E24").Value = 13.28672787389756
$ws.Range("F24").Value = 27.10525974557518
$ws.Range("G24").Value = 3.617072521714422
$ws.Range("J24").Value = 9.594943999741945
$ws.Range("M24").Value = 50.84316916630146
$ws.Range("O24").Value = 21.32900595282469
$ws.Range("B25").Value = 7.899992797380273
$ws.Range("D25").Value = 8.559880643249787
$ws.Range("E25").Value = 13.21004014287652
$ws.Range("F25").Value = 27.65980198590292
$ws.Range("G25").Value = 3.6230351474418
$ws.Range("J25").Value = 9.593809352815924
$ws.Range("M25").Value = 46.74298986400108
$ws.Range("O25").Value = 21.6080511658091

Write-Output "Applied 192 cell updates for Case_1_161 (380 kV)"